$wb = $excel.ActiveWorkbook

# Sheet: Neodymium
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C2").Value = [double]"4.05721109302746E-07"
$ws.Range("D2").Value = 0.1664791307295803
$ws.Range("E2").Value = 0.4065106765944204
$ws.Range("B3").Value = [double]"4.691044125953376E-13"
$ws.Range("C3").Value = 0.001512131801220273
$ws.Range("D3").Value = 0.1058606147983076
$ws.Range("E3").Value = 0.3596305710095618
$ws.Range("B4").Value = [double]"7.321669688613381E-15"
$ws.Range("C4").Value = 0.0003740244660593755
$ws.Range("D4").Value = 0.07931982937362549
$ws.Range("E4").Value = 0.3171595483280041
$ws.Range("C5").Value = [double]"3.409138439013245E-10"
$ws.Range("D5").Value = 0.001223194692503424
$ws.Range("E5").Value = 0.02413018472889736

# Sheet: Dysprosium
$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C2").Value = [double]"4.057211093027732E-07"
$ws.Range("D2").Value = 0.1664791307295915
$ws.Range("E2").Value = 0.4065106765944476
$ws.Range("B3").Value = [double]"4.691044125953693E-13"
$ws.Range("C3").Value = 0.001512131801220375
$ws.Range("D3").Value = 0.1058606147983148
$ws.Range("E3").Value = 0.3596305710095858
$ws.Range("B4").Value = [double]"7.321669688613872E-15"
$ws.Range("C4").Value = 0.0003740244660594007
$ws.Range("D4").Value = 0.07931982937363082
$ws.Range("E4").Value = 0.3171595483280256
$ws.Range("C5").Value = [double]"3.409138439013518E-10"
$ws.Range("D5").Value = 0.001223194692503522
$ws.Range("E5").Value = 0.02413018472889929

# Sheet: Copper
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"6.026889418659203E-07"
$ws.Range("C2").Value = 0.0002786503724174327
$ws.Range("D2").Value = 0.2533655906882806
$ws.Range("E2").Value = 0.6996950518224571
$ws.Range("B3").Value = [double]"4.096143759209795E-06"
$ws.Range("C3").Value = 0.005525867807623958
$ws.Range("D3").Value = 0.15171349876109
$ws.Range("E3").Value = 0.4843334564770065
$ws.Range("B4").Value = [double]"1.21620817902831E-05"
$ws.Range("C4").Value = 0.0007198594464809803
$ws.Range("D4").Value = 0.09476617436401481
$ws.Range("E4").Value = 0.4657667777509991
$ws.Range("B5").Value = [double]"3.819645484553135E-06"
$ws.Range("C5").Value = 0.001678931105475365
$ws.Range("D5").Value = 0.1708471635452753
$ws.Range("E5").Value = 0.4717081075110235

# Sheet: Raw silicon
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"5.031589966768445E-06"
$ws.Range("C2").Value = 0.0003795813608333849
$ws.Range("D2").Value = 0.2641707319413292
$ws.Range("E2").Value = 1.383839401004388
$ws.Range("B3").Value = [double]"5.367099139082082E-06"
$ws.Range("C3").Value = 0.001593037759322987
$ws.Range("D3").Value = 0.1524988814736611
$ws.Range("E3").Value = 0.6174302220107619
$ws.Range("B4").Value = [double]"3.443534520072915E-05"
$ws.Range("C4").Value = 0.000412080869035601
$ws.Range("D4").Value = 0.1111119351128833
$ws.Range("E4").Value = 0.6154403591514438
$ws.Range("B5").Value = [double]"1.848515343462404E-05"
$ws.Range("C5").Value = 0.0005196111533190503
$ws.Range("D5").Value = 0.2361937477439512
$ws.Range("E5").Value = 0.9816170507989759